# Split the single long Bibliografia run into one <w:t> segment per
# reference, joined by manual line breaks (<w:br/>), matching the target
# edit. The visible characters are unchanged -- only line breaks are
# inserted at the reference boundaries (two breaks where a blank line
# separates sections).
$d = $word.ActiveDocument

$range = $d.Content
$find = $range.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$searchText = "Slack, N. & Lewis, M. Operations strategy. Pearson Education, 2020.Slack, N. & Lewis, M. Estratégia de operações. Bookman, 2009.Bibliografia complementarCARVALHO, M. M.; LAURINDO, F.J.B. Estratégia Competitiva: dos conceitos à implementação. São Paulo: Editora Atlas, 2007.MONTGOMERY, C. A. O Estrategista. Editora: Sextante, 2012.BRYNJOLFSSON, E.; HU, Y. J.; SMITH, M. D. From Niches to Riches: Anatomy of the Long Tail. MIT Sloan Management Review, 47, n4, Summer 2006.Dissertação (Engenharia de Produção) - Universidade de São Paulo.2011.CHESBROUGH, H.W.; TEECE, D.J. When is Virtual Virtuous? Harvard Business Review, v.74, n.1, p.65-73, May/June, 1996.CHRISTENSEN, C.M. Making Strategy: Learning by doing. Harvard Business Review, p.141-156, Nov./Dec. 1997.CHRISTENSEN, C.M. The past and future of competitive advantage. MIT Sloan Management Review, v.42, n.2, p105-109, Winter 2001.FLEURY, A.C.C.; FLEURY, M.T.L. Estratégias Empresariais e Formação de Competências: Um quebra-cabeça caleidoscópico da indústria brasileira. São Paulo: Ed. Atlas, 2A ed., 2000.FURRER, O.; SUDHARSHAN, D.; THOMAS, H.; ALEXANDRE, M. T. Resource configurations, generic strategies, and firm performance: Exploring the parallels between resource-based and competitive strategy theories in a new industry. Journal of Strategy and Management,Vol. 1 No. 1, pp. 15-40, 2008.HAMEL, G. Strategic as Revolution. Harvard Business Review, p.69-82, Jul.-Aug., 1996.HAMEL, G; PRAHALAD, C.K. Strategic Intent. Harvard Business Review, p.63-76, May-June, 1989.HENDERSON, B.D. The origin of strategy. Harvard Business Review, Nov/Dec.1989.HURST, D.K. Strategy. Strategy+Business, i. 25, 4th quarter, 2001.KISSIMOTO, K.O. A Influência da Tecnologia da Informação na Estratégia de Personalização nas EmpresasBrasileiras. Dissertação (Engenharia de Produção) - Universidade de São Paulo. 2011.LUEHRMAN, T.A.. ‘Strategy as a Portfolio of Real Options’. Harvard Business Review, p.89-99, Sep./Oct. 1998.MARKIDES C. C. A Dynamic view of strategy. Sloan Management Review, Spring, p., 1999.MAUBORGNE, R.; KIM, W. C. A Estratégia do Oceano Azul: Como Criar Novos Mercados e Tornar a Concorrência Irrelevante. Editora Campus, 258p., 2005.MINTZBERG, H.; LAMPEL, J. Reflecting on the Strategy Process. Sloan management Review, p. 83-94, Spring, 1999.MINTZBERG, H.; AHLSTRAND, B.; LAMPEL, J. Safári de Estratégia, Porto Alegre, Bookman, 299p., 2000.PORTER, M.E. Estratégia Competitiva: Técnicas para Análise de Indústrias e da Concorrência. Editora Campus, 1996, 362p.PORTER, M.E. The Five Competitive Forces that Shape Strategy. Harvard Business Review, V. 86, I1, p78-93, Jan2008.PORTER, M.E. What is Strategy? Harvard Business Review, p.61-78, Nov-Dec, 1996PRAHALAD, C. K.; KRISHNAN, M. S. A Nova Era da Inovação: A Inovação Focada no Relacionamento com o Cliente. Editora Campus, 256 pg., 2008.PRAHALAD, C.K; HAMEL, G. The Core Competence of the Corporation., p.79-91, May-June, 1990.VALIKANGAS, L.; GIBBERT, M. Boundary-Setting Strategies for Escaping Innovation Traps. MIT"
$replaceText = "Slack, N. & Lewis, M. Operations strategy. Pearson Education, 2020.^lSlack, N. & Lewis, M. Estratégia de operações. Bookman, 2009.^l^lBibliografia complementar^l^lCARVALHO, M. M.; LAURINDO, F.J.B. Estratégia Competitiva: dos conceitos à implementação. São Paulo: Editora Atlas, 2007.^lMONTGOMERY, C. A. O Estrategista. Editora: Sextante, 2012.^lBRYNJOLFSSON, E.; HU, Y. J.; SMITH, M. D. From Niches to Riches: Anatomy of the Long Tail. MIT Sloan Management Review, 47, n4, Summer 2006.^lDissertação (Engenharia de Produção) - Universidade de São Paulo.2011.^lCHESBROUGH, H.W.; TEECE, D.J. When is Virtual Virtuous? Harvard Business Review, v.74, n.1, p.65-73, May/June, 1996.^lCHRISTENSEN, C.M. Making Strategy: Learning by doing. Harvard Business Review, p.141-156, Nov./Dec. 1997.^lCHRISTENSEN, C.M. The past and future of competitive advantage. MIT Sloan Management Review, v.42, n.2, p105-109, Winter 2001.^lFLEURY, A.C.C.; FLEURY, M.T.L. Estratégias Empresariais e Formação de Competências: Um quebra-cabeça caleidoscópico da indústria brasileira. São Paulo: Ed. Atlas, 2A ed., 2000.^lFURRER, O.; SUDHARSHAN, D.; THOMAS, H.; ALEXANDRE, M. T. Resource configurations, generic strategies, and firm performance: Exploring the parallels between resource-based and competitive strategy theories in a new industry. Journal of Strategy and Management,Vol. 1 No. 1, pp. 15-40, 2008.^lHAMEL, G. Strategic as Revolution. Harvard Business Review, p.69-82, Jul.-Aug., 1996.^lHAMEL, G; PRAHALAD, C.K. Strategic Intent. Harvard Business Review, p.63-76, May-June, 1989.^lHENDERSON, B.D. The origin of strategy. Harvard Business Review, Nov/Dec.1989.^lHURST, D.K. Strategy. Strategy+Business, i. 25, 4th quarter, 2001.^lKISSIMOTO, K.O. A Influência da Tecnologia da Informação na Estratégia de Personalização nas Empresas^lBrasileiras. Dissertação (Engenharia de Produção) - Universidade de São Paulo. 2011.^lLUEHRMAN, T.A.. ‘Strategy as a Portfolio of Real Options’. Harvard Business Review, p.89-99, Sep./Oct. 1998.^lMARKIDES C. C. A Dynamic view of strategy. Sloan Management Review, Spring, p., 1999.^lMAUBORGNE, R.; KIM, W. C. A Estratégia do Oceano Azul: Como Criar Novos Mercados e Tornar a Concorrência Irrelevante. Editora Campus, 258p., 2005.^lMINTZBERG, H.; LAMPEL, J. Reflecting on the Strategy Process. Sloan management Review, p. 83-94, Spring, 1999.^lMINTZBERG, H.; AHLSTRAND, B.; LAMPEL, J. Safári de Estratégia, Porto Alegre, Bookman, 299p., 2000.^lPORTER, M.E. Estratégia Competitiva: Técnicas para Análise de Indústrias e da Concorrência. Editora Campus, 1996, 362p.^lPORTER, M.E. The Five Competitive Forces that Shape Strategy. Harvard Business Review, V. 86, I1, p78-93, Jan2008.^lPORTER, M.E. What is Strategy? Harvard Business Review, p.61-78, Nov-Dec, 1996^lPRAHALAD, C. K.; KRISHNAN, M. S. A Nova Era da Inovação: A Inovação Focada no Relacionamento com o Cliente. Editora Campus, 256 pg., 2008.^lPRAHALAD, C.K; HAMEL, G. The Core Competence of the Corporation., p.79-91, May-June, 1990.^lVALIKANGAS, L.; GIBBERT, M. Boundary-Setting Strategies for Escaping Innovation Traps. MIT"

$result = $find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)

if (-not $result) {
    throw "Find/Replace on the Bibliografia paragraph did not match; document may have changed."
}

Write-Output "replaced=$result"
